# Auto-generated Excel COM-interop script applying the scheduled-runner
# price/profit update described by the commit diff.
# For each affected row, the H-N "live market" columns are refreshed with
# newly fetched values; a few cells are newly populated or cleared entirely
# to mirror rows whose computed fields became zero/blank (or vice versa).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4072.6365
$ws.Range("J64").Value = 5000
$ws.Range("L64").Value = 5000
$ws.Range("N64").Value = -5496

$ws.Range("H67").Value = 4072.6365
$ws.Range("J67").Value = 5000
$ws.Range("L67").Value = 5000
$ws.Range("N67").Value = -6716

$ws.Range("H113").Value = 3284
$ws.Range("I113").Value = 3079.6
$ws.Range("K113").Value = 3079.6
$ws.Range("M113").Value = 174.4000000000001

$ws.Range("H116").Value = 9185.143
$ws.Range("J116").Value = 6665.6665
$ws.Range("L116").Value = 6665.6665
$ws.Range("N116").Value = -13549.6665

$ws.Range("H131").Value = 6507.278
$ws.Range("I131").Value = 1597.5
$ws.Range("K131").Value = 4792.5
$ws.Range("M131").Value = 247.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 11247635
$ws.Range("I2").Value = 14287191
$ws.Range("K2").Value = 14287191
$ws.Range("M2").Value = -14287078

$ws.Range("H32").Value = 9346.923000000001
$ws.Range("I32").Value = 8003.048
$ws.Range("J32").Value = 14991.2
$ws.Range("K32").Value = 8003.048
$ws.Range("L32").Value = 14991.2
$ws.Range("M32").Value = -7716.048
$ws.Range("N32").Value = -15565.2

$ws.Range("H45").Value = 13808.6
$ws.Range("I45").Value = 16985.75
$ws.Range("J45").Value = 1100
$ws.Range("K45").Value = 16985.75
$ws.Range("L45").Value = 1100
$ws.Range("M45").Value = -16608.75
$ws.Range("N45").Value = -1854

$ws.Range("H61").Value = 8253.333000000001
$ws.Range("I61").Value = 8870.1875
$ws.Range("J61").Value = 6279.4
$ws.Range("K61").Value = 8870.1875
$ws.Range("L61").Value = 6279.4
$ws.Range("M61").Value = -8658.1875
$ws.Range("N61").Value = -6703.4

$ws.Range("H92").Value = 9329.666999999999
$ws.Range("J92").Value = 9329.666999999999
$ws.Range("L92").Value = 9329.666999999999
$ws.Range("N92").Value = -14321.667

$ws.Range("H97").Value = 21281054
$ws.Range("I97").Value = 31254608
$ws.Range("J97").Value = 4138.8
$ws.Range("K97").Value = 31254608
$ws.Range("L97").Value = 4138.8
$ws.Range("M97").Value = -31254112
$ws.Range("N97").Value = -5130.8

$ws.Range("H110").Value = 1916.4166
$ws.Range("I110").Value = 1827.091
$ws.Range("J110").Value = 2899
$ws.Range("K110").Value = 1827.091
$ws.Range("L110").Value = 2899
$ws.Range("M110").Value = 217.9090000000001
$ws.Range("N110").Value = -6989

$ws.Range("H116").Value = 11247635
$ws.Range("I116").Value = 14287191
$ws.Range("K116").Value = 14287191
$ws.Range("M116").Value = -14284897

$ws.Range("H132").Value = 3241.6
$ws.Range("J132").Value = 3711.4285
$ws.Range("L132").Value = 11134.2855
$ws.Range("N132").Value = -16194.2855

$ws.Range("H136").Value = 8253.333000000001
$ws.Range("I136").Value = 8870.1875
$ws.Range("J136").Value = 6279.4
$ws.Range("K136").Value = 26610.5625
$ws.Range("L136").Value = 18838.2
$ws.Range("M136").Value = -24060.5625
$ws.Range("N136").Value = -23938.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 11247635
$ws.Range("I3").Value = 14287191
$ws.Range("K3").Value = 14287191
$ws.Range("M3").Value = -14287077

$ws.Range("H80").Value = 1667574.1
$ws.Range("J80").Value = 2381855
$ws.Range("L80").Value = 2381855
$ws.Range("N80").Value = -2383851

$ws.Range("H83").Value = 1667574.1
$ws.Range("J83").Value = 2381855
$ws.Range("L83").Value = 11909275
$ws.Range("N83").Value = -11919259

$ws.Range("H107").Value = 4930.579
$ws.Range("I107").Value = 4748.385
$ws.Range("K107").Value = 4748.385
$ws.Range("M107").Value = -2828.385

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2422.3208
$ws.Range("I31").Value = 2167.6743
$ws.Range("K31").Value = 2167.6743
$ws.Range("M31").Value = -1872.6743

$ws.Range("H34").Value = 2422.3208
$ws.Range("I34").Value = 2167.6743
$ws.Range("K34").Value = 2167.6743
$ws.Range("M34").Value = -1965.6743

$ws.Range("H59").Value = 58775.5
$ws.Range("J59").Value = 78332.664
$ws.Range("L59").Value = 78332.664
$ws.Range("N59").Value = -80622.664

$ws.Range("H138").Value = 118863.78
$ws.Range("J138").Value = 118863.78
$ws.Range("L138").Value = 118863.78
$ws.Range("N138").Value = -129143.78

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 2074.3333
$ws.Range("I23").Value = 647
$ws.Range("J23").Value = 2788
$ws.Range("K23").Value = 1941
$ws.Range("L23").Value = 8364
$ws.Range("M23").Value = -1706
$ws.Range("N23").Value = -8834

$ws.Range("H138").Value = 3802.6667
$ws.Range("I138").Value = 3802.6667
$ws.Range("K138").Value = 11408.0001
$ws.Range("M138").Value = -6268.000100000001

$ws.Range("H139").Value = 3614.2856
$ws.Range("I139").Value = 1493.5714
$ws.Range("K139").Value = 4480.7142
$ws.Range("M139").Value = 659.2857999999997

$ws.Range("H140").Value = 1420.0555
$ws.Range("I140").Value = 1420.0555
$ws.Range("K140").Value = 4260.166499999999
$ws.Range("M140").Value = 919.8335000000006

$ws.Range("H141").Value = 26249.25
$ws.Range("I141").Value = 26249.25
$ws.Range("K141").Value = 78747.75
$ws.Range("M141").Value = -73567.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 2730.8
$ws.Range("I31").Value = 923.1111
$ws.Range("K31").Value = 923.1111
$ws.Range("M31").Value = -631.1111

$ws.Range("H33").Value = 17457
$ws.Range("I33").Value = 14995
$ws.Range("J33").Value = 19919
$ws.Range("K33").Value = 14995
$ws.Range("L33").Value = 19919
$ws.Range("M33").Value = -14743
$ws.Range("N33").Value = -20423

$ws.Range("H35").Value = 23475
$ws.Range("I35").Value = 23475
$ws.Range("K35").Value = 23475
$ws.Range("M35").Value = -23177

$ws.Range("H36").Value = 12500
$ws.Range("I36").Value = 12500
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 12500
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -12015
$ws.Range("N36").ClearContents()

$ws.Range("H37").Value = 2730.8
$ws.Range("I37").Value = 923.1111
$ws.Range("K37").Value = 923.1111
$ws.Range("M37").Value = -646.1111

$ws.Range("H38").Value = 38000
$ws.Range("J38").Value = 38000
$ws.Range("L38").Value = 38000
$ws.Range("N38").Value = -38926

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H41").Value = 20799.857
$ws.Range("I41").Value = 23766.5
$ws.Range("J41").Value = 3000
$ws.Range("K41").Value = 23766.5
$ws.Range("L41").Value = 3000
$ws.Range("M41").Value = -23411.5
$ws.Range("N41").Value = -3710

$ws.Range("H43").Value = 20000
$ws.Range("I43").Value = 20000
$ws.Range("K43").Value = 20000
$ws.Range("M43").Value = -19849

$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("N59").ClearContents()

$ws.Range("H92").Value = 9250.333000000001
$ws.Range("J92").Value = 9250.333000000001
$ws.Range("L92").Value = 9250.333000000001
$ws.Range("N92").Value = -12994.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1374.4166
$ws.Range("J93").Value = 2399.5
$ws.Range("L93").Value = 2399.5
$ws.Range("N93").Value = -4895.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 12163.8
$ws.Range("I23").Value = 2773.3333
$ws.Range("K23").Value = 2773.3333
$ws.Range("M23").Value = -2544.3333

$ws.Range("H136").Value = 2439.375
$ws.Range("I136").Value = 1706.8334
$ws.Range("K136").Value = 5120.5002
$ws.Range("M136").Value = -2570.5002
